# Auto-generated script to update Sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# per scheduled market-price refresh. Plain numeric cell values; no formulas.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2464.5
$ws.Range("I43").Value = 1988
$ws.Range("J43").Value = 2559.8
$ws.Range("K43").Value = 1988
$ws.Range("L43").Value = 2559.8
$ws.Range("M43").Value = -1919
$ws.Range("N43").Value = -2697.8
$ws.Range("H64").Value = 4987
$ws.Range("J64").Value = 4988.3335
$ws.Range("L64").Value = 4988.3335
$ws.Range("N64").Value = -5484.3335
$ws.Range("H67").Value = 4987
$ws.Range("J67").Value = 4988.3335
$ws.Range("L67").Value = 4988.3335
$ws.Range("N67").Value = -6704.3335
$ws.Range("H125").Value = 1344.25
$ws.Range("I125").Value = 1012.5
$ws.Range("J125").Value = 1454.8334
$ws.Range("K125").Value = 9112.5
$ws.Range("L125").Value = 13093.5006
$ws.Range("M125").Value = -6652.5
$ws.Range("N125").Value = -18013.5006
$ws.Range("H131").Value = 3632.75
$ws.Range("I131").Value = 1294.5714
$ws.Range("K131").Value = 3883.7142
$ws.Range("M131").Value = 1156.2858
$ws.Range("H137").Value = 2888
$ws.Range("I137").Value = 1090
$ws.Range("J137").Value = 3287.5557
$ws.Range("K137").Value = 3270
$ws.Range("L137").Value = 9862.667099999999
$ws.Range("M137").Value = -720
$ws.Range("N137").Value = -14962.6671
$ws.Range("H141").Value = 7151.5713
$ws.Range("I141").Value = 7151.5713
$ws.Range("K141").Value = 21454.7139
$ws.Range("M141").Value = -16274.7139

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2536.9092
$ws.Range("I61").Value = 2536.9092
$ws.Range("K61").Value = 2536.9092
$ws.Range("M61").Value = -2324.9092
$ws.Range("H63").Value = 2341
$ws.Range("I63").Value = 2341
$ws.Range("K63").Value = 2341
$ws.Range("M63").Value = -1655
$ws.Range("H66").Value = 2341
$ws.Range("I66").Value = 2341
$ws.Range("K66").Value = 11705
$ws.Range("M66").Value = -8273
$ws.Range("H74").Value = 59700.117
$ws.Range("I74").Value = 59700.117
$ws.Range("K74").Value = 59700.117
$ws.Range("M74").Value = -58826.117
$ws.Range("H77").Value = 59700.117
$ws.Range("I77").Value = 59700.117
$ws.Range("K77").Value = 298500.585
$ws.Range("M77").Value = -294132.585
$ws.Range("H122").Value = 1919.9231
$ws.Range("I122").Value = 1473.2222
$ws.Range("J122").Value = 2925
$ws.Range("K122").Value = 4419.6666
$ws.Range("L122").Value = 8775
$ws.Range("M122").Value = -1969.6666
$ws.Range("N122").Value = -13675
$ws.Range("H132").Value = 68616.92999999999
$ws.Range("I132").Value = 73303.86
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 219911.58
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -217381.58
$ws.Range("N132").Value = -14060
$ws.Range("H136").Value = 2536.9092
$ws.Range("I136").Value = 2536.9092
$ws.Range("K136").Value = 7610.7276
$ws.Range("M136").Value = -5060.7276

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 456.83334
$ws.Range("J80").Value = 339.2
$ws.Range("L80").Value = 339.2
$ws.Range("N80").Value = -2335.2
$ws.Range("H82").Value = 9396.556
$ws.Range("I82").Value = 9396.556
$ws.Range("K82").Value = 9396.556
$ws.Range("M82").Value = -9013.556
$ws.Range("H83").Value = 456.83334
$ws.Range("J83").Value = 339.2
$ws.Range("L83").Value = 1696
$ws.Range("N83").Value = -11680
$ws.Range("H85").Value = 9396.556
$ws.Range("I85").Value = 9396.556
$ws.Range("K85").Value = 9396.556
$ws.Range("M85").Value = -8070.556
$ws.Range("H134").Value = 2785.5
$ws.Range("I134").Value = 2219.8
$ws.Range("K134").Value = 6659.400000000001
$ws.Range("M134").Value = -4124.400000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 649.55
$ws.Range("I16").Value = 649.55
$ws.Range("K16").Value = 649.55
$ws.Range("M16").Value = -362.55
$ws.Range("H22").Value = 1348.7368
$ws.Range("I22").Value = 289.91666
$ws.Range("K22").Value = 289.91666
$ws.Range("M22").Value = 60.08334000000002
$ws.Range("H58").Value = 92727.45
$ws.Range("I58").Value = 101598.8
$ws.Range("K58").Value = 101598.8
$ws.Range("M58").Value = -101395.8
$ws.Range("H113").Value = 649.55
$ws.Range("I113").Value = 649.55
$ws.Range("K113").Value = 649.55
$ws.Range("M113").Value = 1520.45
$ws.Range("H122").Value = 3019.8
$ws.Range("I122").Value = 2999
$ws.Range("J122").Value = 3025
$ws.Range("K122").Value = 8997
$ws.Range("L122").Value = 9075
$ws.Range("M122").Value = -6547
$ws.Range("N122").Value = -13975
$ws.Range("H132").Value = 2363.2
$ws.Range("I132").Value = 2363.2
$ws.Range("K132").Value = 7089.599999999999
$ws.Range("M132").Value = -4559.599999999999
$ws.Range("H134").Value = 94642
$ws.Range("I134").Value = 168205.67
$ws.Range("K134").Value = 504617.01
$ws.Range("M134").Value = -502082.01
$ws.Range("H136").Value = 92727.45
$ws.Range("I136").Value = 101598.8
$ws.Range("K136").Value = 304796.4
$ws.Range("M136").Value = -302246.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1658.4117
$ws.Range("I34").Value = 1413.9286
$ws.Range("J34").Value = 2799.3333
$ws.Range("K34").Value = 4241.7858
$ws.Range("L34").Value = 8397.999899999999
$ws.Range("M34").Value = -4157.7858
$ws.Range("N34").Value = -8565.999899999999
$ws.Range("H55").Value = 8080
$ws.Range("I55").Value = 4244
$ws.Range("J55").Value = 9998
$ws.Range("K55").Value = 12732
$ws.Range("L55").Value = 29994
$ws.Range("M55").Value = -12555
$ws.Range("N55").Value = -30348
$ws.Range("H80").Value = 5257.6665
$ws.Range("J80").Value = 5257.6665
$ws.Range("L80").Value = 15772.9995
$ws.Range("N80").Value = -17644.9995
$ws.Range("H83").Value = 5257.6665
$ws.Range("J83").Value = 5257.6665
$ws.Range("L83").Value = 47318.9985
$ws.Range("N83").Value = -56678.9985

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H29").Value = 954.3333
$ws.Range("I29").Value = 1006.5
$ws.Range("J29").Value = 850
$ws.Range("K29").Value = 1006.5
$ws.Range("L29").Value = 850
$ws.Range("M29").Value = -716.5
$ws.Range("N29").Value = -1430
$ws.Range("H80").Value = 3258.875
$ws.Range("I80").Value = 2595.2222
$ws.Range("J80").Value = 3657.0667
$ws.Range("K80").Value = 2595.2222
$ws.Range("L80").Value = 3657.0667
$ws.Range("M80").Value = -1597.2222
$ws.Range("N80").Value = -5653.066699999999
$ws.Range("H83").Value = 3258.875
$ws.Range("I83").Value = 2595.2222
$ws.Range("J83").Value = 3657.0667
$ws.Range("K83").Value = 12976.111
$ws.Range("L83").Value = 18285.3335
$ws.Range("M83").Value = -7984.111000000001
$ws.Range("N83").Value = -28269.3335
$ws.Range("H132").Value = 28254.922
$ws.Range("I132").Value = 33885.676
$ws.Range("K132").Value = 101657.028
$ws.Range("M132").Value = -99127.02799999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4301.88
$ws.Range("I7").Value = 3406
$ws.Range("K7").Value = 3406
$ws.Range("M7").Value = -3294
$ws.Range("H69").Value = 63000
$ws.Range("I69").Value = 37000
$ws.Range("J69").Value = 89000
$ws.Range("K69").Value = 37000
$ws.Range("L69").Value = 89000
$ws.Range("M69").Value = -36189
$ws.Range("N69").Value = -90622
$ws.Range("H72").Value = 63000
$ws.Range("I72").Value = 37000
$ws.Range("J72").Value = 89000
$ws.Range("K72").Value = 111000
$ws.Range("L72").Value = 267000
$ws.Range("M72").Value = -106944
$ws.Range("N72").Value = -275112
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("H126").Value = 4301.88
$ws.Range("I126").Value = 3406
$ws.Range("K126").Value = 10218
$ws.Range("M126").Value = -7748
$ws.Range("H136").Value = 8358.5
$ws.Range("I136").Value = 7811.3335
$ws.Range("K136").Value = 23434.0005
$ws.Range("M136").Value = -20884.0005
$ws.Range("N106").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 40000
$ws.Range("I75").Value = 40000
$ws.Range("K75").Value = 40000
$ws.Range("M75").Value = -39064
$ws.Range("H78").Value = 40000
$ws.Range("I78").Value = 40000
$ws.Range("K78").Value = 120000
$ws.Range("M78").Value = -115320
$ws.Range("H126").Value = 41753.57
$ws.Range("J126").Value = 6412.7144
$ws.Range("L126").Value = 19238.1432
$ws.Range("N126").Value = -24178.1432
$ws.Range("H135").Value = 68810
$ws.Range("J135").Value = 68810
$ws.Range("L135").Value = 68810
$ws.Range("N135").Value = -78950

